# Add a new "caseDetailQuery" column (ICDC/CTDC case level detail switch +
# the Neo4j query used to read the case detail table) in between the
# existing "StatQuery" (col B) and "dbExcel" (old col C) columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new column at C - this shifts the old C (dbExcel) and
# D (WebExcel) columns one slot to the right (-> D and E) and copies the
# left-hand column's formatting (width / wrap-text style) onto the new
# column, same as Excel's normal "Insert Sheet Columns" behaviour.
$ws.Columns("C").Insert()

# Match column C's width to the (identical) width already used by columns
# A and B.
$ws.Columns("C").ColumnWidth = 75

# New header (row 1) + query text (row 2) for the case-detail query.
$ws.Range("C1").Value = "caseDetailQuery"
$ws.Range("C2").Value = "MATCH (f:file)-[*]->(c:case) WITH DISTINCT(f) AS f, c MATCH (f)-->(parent) WHERE c.case_id IN ['NCATS-COP01CCB010072'] RETURN f.file_name AS ``File Name`` ,f.file_type AS ``File Type``,head(labels(parent)) AS ``Association``, f.file_description AS ``Description``,f.file_format AS Format,((f.file_size)/1024) AS Size"

# Leave the cursor/selection where the author left it after the edit.
$ws.Range("B4").Select()
